$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace the "Upon receiving the plaintext_dictionary..." paragraph with
#    four paragraphs: "Encryption Scheme:" heading, the new encryption
#    methodology paragraph, "Decryption Scheme:" heading, and the (slightly
#    reworded-in-markup) original "Upon receiving the plaintext_dictionary..."
#    paragraph.
# ---------------------------------------------------------------------------

$wOpenXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-WordXmlFragment([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" +
        '<w:wordDocument ' + $wOpenXmlNs + '><w:body>' + $bodyInner + '</w:body></w:wordDocument>'
}

# Find the target paragraph by its text (it is currently Paragraph index 2).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Upon receiving the plaintext_dictionary*") {
        $target = $p
        break
    }
}

$targetRange = $target.Range
$insertionPoint = $targetRange.Duplicate
$insertionPoint.Collapse(1)  # wdCollapseStart

# Insert three blank paragraph shells before the target paragraph, then fill
# each one with its final (fully-formatted) content via InsertXML - this
# keeps the surrounding document structure untouched while giving us full
# control over run/formatting XML.
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()

# Re-resolve the target paragraph's index now that three blank paragraphs
# precede it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Upon receiving the plaintext_dictionary*") {
        $target = $p
        break
    }
}
$targetIndex = $target.Index

$encryptionHeadingPara = $d.Paragraphs($targetIndex - 3)
$encryptionBodyPara    = $d.Paragraphs($targetIndex - 2)
$decryptionHeadingPara = $d.Paragraphs($targetIndex - 1)
$finalPara             = $d.Paragraphs($targetIndex)

$encryptionHeadingXml = New-WordXmlFragment '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Encryption Scheme:</w:t></w:r></w:p>'
$encryptionHeadingPara.Range.InsertXML($encryptionHeadingXml)

$encryptionBodyInner = @'
<w:p>
<w:r><w:t xml:space="preserve">Given the pseudocode for the encryption scheme, all three of us created our own encryption methods.  Our encryption schemes used the C++ library </w:t></w:r>
<w:r><w:t>&lt;</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>math.h</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>&gt;</w:t></w:r>
<w:r><w:t xml:space="preserve"> which included a uniform distribution random number generator</w:t></w:r>
<w:r><w:t>.  We used the random number generator to create a set of numbers which were converted into the set of alphabetic letters with an additional case to take into account the possibility of generating a blank space as well</w:t></w:r>
<w:r><w:t>.  This same uniform random n umber generator was used for the coin generation algorithm which would spawn either 1 or 0</w:t></w:r>
<w:r><w:t xml:space="preserve">.  When our individual encryption algorithms were capable of converting plaintext messages into ciphertexts, we combined our algorithms to generate a set of ciphertext messages </w:t></w:r>
</w:p>
'@
$encryptionBodyPara.Range.InsertXML((New-WordXmlFragment $encryptionBodyInner))

$decryptionHeadingInner = '<w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Decryption Scheme</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>'
$decryptionHeadingPara.Range.InsertXML((New-WordXmlFragment $decryptionHeadingInner))

$finalParaInner = @'
<w:p>
<w:r><w:t xml:space="preserve">Upon receiving the </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>plaintext_dictionary</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve">, we conducted a letter-based frequency analysis on all five candidate plaintext messages to use as reference for what the ciphertext may hold.  While the analysis is unable to take into account the possibility of which letters are randomly generated and which letters are encoded, there is a distinctly imbalanced distribution of letters as shown in the graph below.  While all five plaintext candidates have very similar distributions, it identifies and allows us to establish a default mapping of each plaintext letter might map to in regards to the ciphertext.</w:t></w:r>
</w:p>
'@
$finalPara.Range.InsertXML((New-WordXmlFragment $finalParaInner))

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker from the Plaintext #4
#    frequency table's "A" cell to the Plaintext #1 frequency table's "A"
#    cell (both are Cell(1,1) of their respective tables).
# ---------------------------------------------------------------------------

$table1 = $d.Tables(1)
$cellA1 = $table1.Cell(1, 1)
$cellA1ParaRange = $cellA1.Range.Paragraphs(1).Range
$cellA1Xml = New-WordXmlFragment '<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>A</w:t></w:r></w:p>'
$cellA1ParaRange.InsertXML($cellA1Xml)

$table4 = $d.Tables(4)
$cellA4 = $table4.Cell(1, 1)
$cellA4ParaRange = $cellA4.Range.Paragraphs(1).Range
$cellA4Xml = New-WordXmlFragment '<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>A</w:t></w:r></w:p>'
$cellA4ParaRange.InsertXML($cellA4Xml)

# ---------------------------------------------------------------------------
# 3) Add a <w:lastRenderedPageBreak/> immediately before the <w:drawing> in
#    the inline picture's run, without disturbing the drawing/image itself.
# ---------------------------------------------------------------------------

$shape = $d.InlineShapes(1)
$pictureParaRange = $shape.Range.Paragraphs(1).Range
$pictureXml = $pictureParaRange.WordOpenXML
$marker = '<w:drawing>'
$markerIndex = $pictureXml.IndexOf($marker)
if ($markerIndex -ge 0) {
    $updatedPictureXml = $pictureXml.Substring(0, $markerIndex) + '<w:lastRenderedPageBreak/>' + $pictureXml.Substring($markerIndex)
    $pictureParaRange.InsertXML($updatedPictureXml)
}

Write-Output "Done."
